$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header D1 "Publicado em" matching style of C1 (bold/border header style)
$ws.Range("D1").Value = "Publicado em"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "Publicado em"
$excel.CutCopyMode = $false

# Set D2:D393 "Publicado em" date values as plain text (not auto-converted to Excel dates)
$rng = $ws.Range("D2:D393")
$rng.NumberFormat = "@"
$values = New-Object 'object[,]' 392,1
$values[0,0] = "2021-08-13"
$values[1,0] = "2022-05-16"
$values[2,0] = "2022-07-20"
$values[3,0] = "2024-07-31"
$values[4,0] = "2025-01-29"
$values[5,0] = "2022-09-14"
$values[6,0] = "2025-01-10"
$values[7,0] = "2024-10-09"
$values[8,0] = "2024-10-16"
$values[9,0] = "2025-03-03"
$values[10,0] = "2024-12-18"
$values[11,0] = "2025-01-07"
$values[12,0] = "2018-08-27"
$values[13,0] = "2018-08-27"
$values[14,0] = "2022-06-08"
$values[15,0] = "2022-08-16"
$values[16,0] = "2020-07-06"
$values[17,0] = "2025-03-07"
$values[18,0] = "2024-12-18"
$values[19,0] = "2024-09-05"
$values[20,0] = "2024-10-14"
$values[21,0] = "2024-11-11"
$values[22,0] = "2025-02-06"
$values[23,0] = "2024-04-30"
$values[24,0] = "2023-12-12"
$values[25,0] = "2024-06-03"
$values[26,0] = "2024-09-30"
$values[27,0] = "2024-11-08"
$values[28,0] = "2024-11-08"
$values[29,0] = "2024-01-26"
$values[30,0] = "2021-10-12"
$values[31,0] = "2020-04-06"
$values[32,0] = "2021-08-09"
$values[33,0] = "2021-08-09"
$values[34,0] = "2024-07-20"
$values[35,0] = "2024-09-30"
$values[36,0] = "2024-09-30"
$values[37,0] = "2024-03-18"
$values[38,0] = "2024-05-30"
$values[39,0] = "2024-05-30"
$values[40,0] = "2024-07-18"
$values[41,0] = "2024-06-19"
$values[42,0] = "2022-10-07"
$values[43,0] = "2024-09-12"
$values[44,0] = "2024-06-28"
$values[45,0] = "2024-10-26"
$values[46,0] = "2024-11-25"
$values[47,0] = "2024-05-20"
$values[48,0] = "2025-01-09"
$values[49,0] = "2023-02-28"
$values[50,0] = "2024-03-08"
$values[51,0] = "2024-12-16"
$values[52,0] = "2025-01-23"
$values[53,0] = "2025-02-06"
$values[54,0] = "2025-01-22"
$values[55,0] = "2024-11-27"
$values[56,0] = "2025-02-03"
$values[57,0] = "2024-07-05"
$values[58,0] = "2024-07-09"
$values[59,0] = "2024-06-28"
$values[60,0] = "2024-06-25"
$values[61,0] = "2024-07-09"
$values[62,0] = "2024-08-31"
$values[63,0] = "2017-04-05"
$values[64,0] = "2017-04-10"
$values[65,0] = "2017-12-20"
$values[66,0] = "2016-05-04"
$values[67,0] = "2016-12-19"
$values[68,0] = "2019-12-04"
$values[69,0] = "2024-04-15"
$values[70,0] = "2024-05-10"
$values[71,0] = "2022-04-28"
$values[72,0] = "2024-03-07"
$values[73,0] = "2024-07-25"
$values[74,0] = "2024-11-07"
$values[75,0] = "2024-04-03"
$values[76,0] = "2024-01-24"
$values[77,0] = "2024-07-01"
$values[78,0] = "2019-02-28"
$values[79,0] = "2019-02-28"
$values[80,0] = "2020-07-22"
$values[81,0] = "2024-09-29"
$values[82,0] = "2024-10-15"
$values[83,0] = "2025-01-14"
$values[84,0] = "2024-04-15"
$values[85,0] = "2024-10-04"
$values[86,0] = "2024-08-14"
$values[87,0] = "2025-01-31"
$values[88,0] = "2024-12-05"
$values[89,0] = "2024-09-16"
$values[90,0] = "2025-01-17"
$values[91,0] = "2025-02-20"
$values[92,0] = "2024-10-10"
$values[93,0] = "2024-10-31"
$values[94,0] = "2025-01-09"
$values[95,0] = "2022-08-12"
$values[96,0] = "2022-08-15"
$values[97,0] = "2022-08-15"
$values[98,0] = "2022-09-16"
$values[99,0] = "2023-02-06"
$values[100,0] = "2024-10-31"
$values[101,0] = "2024-11-16"
$values[102,0] = "2025-03-12"
$values[103,0] = "2023-11-30"
$values[104,0] = "2024-09-17"
$values[105,0] = "2025-02-28"
$values[106,0] = "2024-10-04"
$values[107,0] = "2022-06-06"
$values[108,0] = "2022-10-10"
$values[109,0] = "2022-06-06"
$values[110,0] = "2024-12-03"
$values[111,0] = "2025-02-25"
$values[112,0] = "2024-08-31"
$values[113,0] = "2024-12-02"
$values[114,0] = "2024-12-03"
$values[115,0] = "2025-02-25"
$values[116,0] = "2024-05-06"
$values[117,0] = "2025-03-28"
$values[118,0] = "2024-06-04"
$values[119,0] = "2024-10-18"
$values[120,0] = "2024-10-31"
$values[121,0] = "2024-05-02"
$values[122,0] = "2024-05-24"
$values[123,0] = "2024-06-10"
$values[124,0] = "2024-10-10"
$values[125,0] = "2024-12-06"
$values[126,0] = "2024-09-18"
$values[127,0] = "2024-09-23"
$values[128,0] = "2024-08-13"
$values[129,0] = "2025-01-31"
$values[130,0] = "2024-08-20"
$values[131,0] = "2024-08-05"
$values[132,0] = "2024-08-05"
$values[133,0] = "2024-08-05"
$values[134,0] = "2024-08-05"
$values[135,0] = "2024-08-05"
$values[136,0] = "2024-08-05"
$values[137,0] = "2020-01-23"
$values[138,0] = "2021-03-03"
$values[139,0] = "2023-01-30"
$values[140,0] = "2023-03-20"
$values[141,0] = "2023-12-04"
$values[142,0] = "2022-12-14"
$values[143,0] = "2022-12-14"
$values[144,0] = "2023-11-03"
$values[145,0] = "2023-11-03"
$values[146,0] = "2023-11-24"
$values[147,0] = "2023-11-24"
$values[148,0] = "2024-01-22"
$values[149,0] = "2024-06-13"
$values[150,0] = "2024-08-05"
$values[151,0] = "2024-09-09"
$values[152,0] = "2025-02-26"
$values[153,0] = "2022-03-14"
$values[154,0] = "2023-02-06"
$values[155,0] = "2023-06-21"
$values[156,0] = "2023-07-11"
$values[157,0] = "2023-07-26"
$values[158,0] = "2024-12-19"
$values[159,0] = "2025-03-03"
$values[160,0] = "2025-02-21"
$values[161,0] = "2024-09-19"
$values[162,0] = "2024-09-20"
$values[163,0] = "2024-11-11"
$values[164,0] = "2024-08-30"
$values[165,0] = "2023-11-13"
$values[166,0] = "2023-11-13"
$values[167,0] = "2023-11-13"
$values[168,0] = "2023-11-13"
$values[169,0] = "2023-11-13"
$values[170,0] = "2023-11-13"
$values[171,0] = "2023-11-20"
$values[172,0] = "2023-12-01"
$values[173,0] = "2023-12-01"
$values[174,0] = "2023-09-19"
$values[175,0] = "2025-02-05"
$values[176,0] = "2022-01-10"
$values[177,0] = "2020-07-13"
$values[178,0] = "2024-07-24"
$values[179,0] = "2024-09-30"
$values[180,0] = "2024-09-30"
$values[181,0] = "2024-09-30"
$values[182,0] = "2025-01-22"
$values[183,0] = "2024-12-18"
$values[184,0] = "2024-10-04"
$values[185,0] = "2024-06-13"
$values[186,0] = "2022-02-16"
$values[187,0] = "2022-02-16"
$values[188,0] = "2022-02-16"
$values[189,0] = "2022-09-29"
$values[190,0] = "2022-09-29"
$values[191,0] = "2022-09-29"
$values[192,0] = "2023-05-09"
$values[193,0] = "2023-05-15"
$values[194,0] = "2024-03-07"
$values[195,0] = "2023-11-06"
$values[196,0] = "2023-11-06"
$values[197,0] = "2023-11-06"
$values[198,0] = "2023-11-06"
$values[199,0] = "2023-11-06"
$values[200,0] = "2023-11-13"
$values[201,0] = "2023-12-18"
$values[202,0] = "2023-12-18"
$values[203,0] = "2024-05-14"
$values[204,0] = "2021-01-15"
$values[205,0] = "2021-01-15"
$values[206,0] = "2024-09-04"
$values[207,0] = "2023-04-24"
$values[208,0] = "2023-02-16"
$values[209,0] = "2023-02-16"
$values[210,0] = "2018-07-12"
$values[211,0] = "2018-07-17"
$values[212,0] = "2020-09-23"
$values[213,0] = "2020-09-23"
$values[214,0] = "2024-04-15"
$values[215,0] = "2024-10-07"
$values[216,0] = "2023-10-27"
$values[217,0] = "2024-06-12"
$values[218,0] = "2024-06-14"
$values[219,0] = "2024-06-14"
$values[220,0] = "2024-06-14"
$values[221,0] = "2024-06-14"
$values[222,0] = "2024-06-14"
$values[223,0] = "2024-11-06"
$values[224,0] = "2024-05-16"
$values[225,0] = "2024-07-17"
$values[226,0] = "2023-09-22"
$values[227,0] = "2023-09-28"
$values[228,0] = "2024-10-09"
$values[229,0] = "2024-06-04"
$values[230,0] = "2024-06-04"
$values[231,0] = "2024-06-23"
$values[232,0] = "2024-09-05"
$values[233,0] = "2024-09-05"
$values[234,0] = "2024-03-26"
$values[235,0] = "2024-08-20"
$values[236,0] = "2024-06-05"
$values[237,0] = "2024-06-05"
$values[238,0] = "2024-06-05"
$values[239,0] = "2024-06-05"
$values[240,0] = "2024-09-09"
$values[241,0] = "2024-09-09"
$values[242,0] = "2024-09-10"
$values[243,0] = "2024-09-10"
$values[244,0] = "2024-07-25"
$values[245,0] = "2024-07-25"
$values[246,0] = "2017-08-31"
$values[247,0] = "2022-04-27"
$values[248,0] = "2021-01-01"
$values[249,0] = "2023-05-03"
$values[250,0] = "2023-05-03"
$values[251,0] = "2024-09-17"
$values[252,0] = "2024-09-17"
$values[253,0] = "2024-05-28"
$values[254,0] = "2024-03-26"
$values[255,0] = "2024-04-11"
$values[256,0] = "2023-04-11"
$values[257,0] = "2023-08-23"
$values[258,0] = "2023-09-21"
$values[259,0] = "2021-02-22"
$values[260,0] = "2021-03-17"
$values[261,0] = "2021-03-17"
$values[262,0] = "2021-03-17"
$values[263,0] = "2021-03-17"
$values[264,0] = "2022-08-16"
$values[265,0] = "2022-08-16"
$values[266,0] = "2022-08-16"
$values[267,0] = "2023-11-20"
$values[268,0] = "2023-12-07"
$values[269,0] = "2024-06-08"
$values[270,0] = "2019-11-26"
$values[271,0] = "2016-04-18"
$values[272,0] = "2016-11-23"
$values[273,0] = "2020-06-15"
$values[274,0] = "2024-11-28"
$values[275,0] = "2024-09-16"
$values[276,0] = "2024-10-22"
$values[277,0] = "2024-10-16"
$values[278,0] = "2025-01-10"
$values[279,0] = "2024-07-23"
$values[280,0] = "2024-03-05"
$values[281,0] = "2024-07-10"
$values[282,0] = "2024-04-11"
$values[283,0] = "2024-04-23"
$values[284,0] = "2024-04-23"
$values[285,0] = "2022-07-19"
$values[286,0] = "2022-05-30"
$values[287,0] = "2023-12-15"
$values[288,0] = "2024-06-06"
$values[289,0] = "2024-10-10"
$values[290,0] = "2019-04-25"
$values[291,0] = "2022-04-18"
$values[292,0] = "2024-12-12"
$values[293,0] = "2018-07-31"
$values[294,0] = "2021-06-08"
$values[295,0] = "2021-08-05"
$values[296,0] = "2022-08-24"
$values[297,0] = "2021-10-22"
$values[298,0] = "2022-08-09"
$values[299,0] = "2023-03-23"
$values[300,0] = "2021-08-31"
$values[301,0] = "2019-06-19"
$values[302,0] = "2019-07-31"
$values[303,0] = "2024-01-09"
$values[304,0] = "2019-01-21"
$values[305,0] = "2021-05-20"
$values[306,0] = "2021-05-21"
$values[307,0] = "2021-05-21"
$values[308,0] = "2021-12-01"
$values[309,0] = "2021-12-01"
$values[310,0] = "2021-12-01"
$values[311,0] = "2021-12-14"
$values[312,0] = "2021-12-16"
$values[313,0] = "2021-12-16"
$values[314,0] = "2024-06-20"
$values[315,0] = "2024-07-02"
$values[316,0] = "2024-07-02"
$values[317,0] = "2024-10-15"
$values[318,0] = "2024-10-15"
$values[319,0] = "2022-09-01"
$values[320,0] = "2023-08-18"
$values[321,0] = "2021-06-11"
$values[322,0] = "2017-11-07"
$values[323,0] = "2021-10-19"
$values[324,0] = "2019-02-18"
$values[325,0] = "2021-01-20"
$values[326,0] = "2021-11-08"
$values[327,0] = "2021-04-29"
$values[328,0] = "2024-06-17"
$values[329,0] = "2024-09-13"
$values[330,0] = "2024-11-21"
$values[331,0] = "2024-10-03"
$values[332,0] = "2024-11-26"
$values[333,0] = "2024-02-29"
$values[334,0] = "2024-02-29"
$values[335,0] = "2023-03-22"
$values[336,0] = "2023-05-31"
$values[337,0] = "2023-01-31"
$values[338,0] = "2022-08-01"
$values[339,0] = "2022-08-01"
$values[340,0] = "2023-01-18"
$values[341,0] = "2024-01-12"
$values[342,0] = "2023-11-16"
$values[343,0] = "2023-11-20"
$values[344,0] = "2024-02-14"
$values[345,0] = "2015-09-29"
$values[346,0] = "2025-01-02"
$values[347,0] = "2022-09-19"
$values[348,0] = "2018-10-19"
$values[349,0] = "2016-11-18"
$values[350,0] = "2015-05-18"
$values[351,0] = "2020-09-21"
$values[352,0] = "2021-04-05"
$values[353,0] = "2021-05-25"
$values[354,0] = "2020-11-26"
$values[355,0] = "2021-01-12"
$values[356,0] = "2021-01-12"
$values[357,0] = "2023-08-31"
$values[358,0] = "2019-02-18"
$values[359,0] = "2017-10-19"
$values[360,0] = "2020-09-10"
$values[361,0] = "2020-09-11"
$values[362,0] = "2013-05-29"
$values[363,0] = "2016-11-01"
$values[364,0] = "2021-12-15"
$values[365,0] = "2020-03-12"
$values[366,0] = "2020-05-13"
$values[367,0] = "2021-04-27"
$values[368,0] = "2023-01-17"
$values[369,0] = "2020-03-28"
$values[370,0] = "2018-10-21"
$values[371,0] = "2024-10-29"
$values[372,0] = "2016-08-12"
$values[373,0] = "2022-05-09"
$values[374,0] = "2021-11-05"
$values[375,0] = "2022-03-17"
$values[376,0] = "2021-02-05"
$values[377,0] = "2021-04-26"
$values[378,0] = "2020-11-16"
$values[379,0] = "2022-08-25"
$values[380,0] = "2020-06-23"
$values[381,0] = "2022-01-19"
$values[382,0] = "2022-01-19"
$values[383,0] = "2015-02-18"
$values[384,0] = "2024-08-23"
$values[385,0] = "2024-08-23"
$values[386,0] = "2024-10-01"
$values[387,0] = "2024-10-01"
$values[388,0] = "2022-05-01"
$values[389,0] = "2022-05-01"
$values[390,0] = "2022-05-17"
$values[391,0] = "2022-05-17"
$rng.Value = $values
$rng.Style = "Normal"

$wb.Save()
